# feat: add 2022-Q1 data
#
# The former "总计" (roll-up) sheet is renamed to "2022-Q1" and repopulated
# with the per-fund holding detail for that quarter (same A:H layout used by
# the "2020-Q4" / "2021-Q1" / "2021-Q4" sheets). A brand-new "总计" sheet is
# appended at the end of the workbook holding the refreshed roll-up table
# (a new 2022-Q1 row on top, the previous rows shifted down by one).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$detail = $wb.Worksheets.Item("总计")
$detail.Name = "2022-Q1"
$detail.Cells.Clear()

# Reuse the header / index-column formatting (bold, bordered, centered) from
# the sibling quarter sheet instead of constructing a brand-new style.
$styleSrc.Range("B1:H1").Copy()
$detail.Range("B1:H1").PasteSpecial($xlPasteFormats)
$styleSrc.Range("A2").Copy()
$detail.Range("A2:A4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$detail.Range("B1").Value = "基金代码"
$detail.Range("C1").Value = "基金名称"
$detail.Range("D1").Value = "基金规模"
$detail.Range("E1").Value = "股票总仓位"
$detail.Range("F1").Value = "仓位占比"
$detail.Range("G1").Value = "持有市值(亿元)"
$detail.Range("H1").Value = "仓位排名"

# Fund code / scale / position figures are stored as plain text in the
# source data (e.g. leading-zero codes like "000586"), so force text before
# writing the values.
$detail.Range("B2:B4").NumberFormat = "@"
$detail.Range("D2:G4").NumberFormat = "@"

$detailRows = @(
    @(0, "000586", "景顺长城中小板创业板精选股票", "2.42", "94.15", "6.12", "0.1481", 7),
    @(1, "010706", "景顺长城景骊成长混合型证券投资基金", "1.13", "93.50", "5.85", "0.0661", 6),
    @(2, "260115", "景顺长城中小盘混合", "0.96", "94.00", "5.14", "0.0493", 9)
)

foreach ($r in $detailRows) {
    $row = $r[0] + 2
    $detail.Cells.Item($row, 1).Value = $r[0]
    $detail.Cells.Item($row, 2).Value = $r[1]
    $detail.Cells.Item($row, 3).Value = $r[2]
    $detail.Cells.Item($row, 4).Value = $r[3]
    $detail.Cells.Item($row, 5).Value = $r[4]
    $detail.Cells.Item($row, 6).Value = $r[5]
    $detail.Cells.Item($row, 7).Value = $r[6]
    $detail.Cells.Item($row, 8).Value = $r[7]
}

# ---------------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet (after all quarter sheets) with the
# updated roll-up table.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$summary.Name = "总计"

$styleSrc.Range("B1:D1").Copy()
$summary.Range("B1:D1").PasteSpecial($xlPasteFormats)
$styleSrc.Range("A2").Copy()
$summary.Range("A2:A5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$summary.Range("B1").Value = "日期"
$summary.Range("C1").Value = "持有数量(只)"
$summary.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @(0, "2022-Q1", 3, 0.26),
    @(1, "2021-Q4", 2, 0.26),
    @(2, "2021-Q1", 3, 1.12),
    @(3, "2020-Q4", 2, 1.01)
)

foreach ($r in $summaryRows) {
    $row = $r[0] + 2
    $summary.Cells.Item($row, 1).Value = $r[0]
    $summary.Cells.Item($row, 2).Value = $r[1]
    $summary.Cells.Item($row, 3).Value = $r[2]
    $summary.Cells.Item($row, 4).Value = $r[3]
}

# Leave the originally-active first sheet selected, matching the
# pre-edit workbook (the edit itself doesn't touch sheet selection).
$wb.Worksheets.Item(1).Activate()
